$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F10").Value = 5336
$ws1.Range("F16").Value = 191

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F10").Value = 5336
$ws4.Range("F16").Value = 191
